$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3.05
$ws.Range("N2").Value = 2.98
$ws.Range("P2").Value = 1.65
$ws.Range("F3").Value = 1.47
$ws.Range("G4").Value = 1.36
$ws.Range("H4").Value = 1.09
$ws.Range("N4").Value = 3.7
$ws.Range("P4").Value = 2.04
$ws.Range("W4").Value = 3.7
$ws.Range("H5").Value = 6.2
$ws.Range("M5").Value = 1.06
$ws.Range("T5").Value = 2.12
$ws.Range("U5").Value = 1.68
$ws.Range("U6").Value = 2.3
$ws.Range("AA6").Value = 46
$ws.Range("AF6").Value = 25
$ws.Range("AI6").Value = 44
$ws.Range("AK6").Value = 34
$ws.Range("I7").Value = 1.28
$ws.Range("M7").Value = 1.02
$ws.Range("F8").Value = 2.48
$ws.Range("J8").Value = 3
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 2.74
$ws.Range("Q8").Value = 1.98
$ws.Range("F9").Value = 1.86
$ws.Range("G9").Value = 1.98
$ws.Range("H9").Value = 4.4
$ws.Range("I9").Value = 5
$ws.Range("K9").Value = 3.9
$ws.Range("O9").Value = 1.32
$ws.Range("P9").Value = 1.86
$ws.Range("R9").Value = 1.33
$ws.Range("V9").Value = 1.25
$ws.Range("Z9").Value = 38
$ws.Range("AA9").Value = 140
$ws.Range("AC9").Value = 8.800000000000001
$ws.Range("AD9").Value = 19.5
$ws.Range("AK9").Value = 22
$ws.Range("AL9").Value = 40
$ws.Range("AM9").Value = 140
$ws.Range("AO9").Value = 95
$ws.Range("F10").Value = 2.4
$ws.Range("G10").Value = 2.62
$ws.Range("H10").Value = 2.78
$ws.Range("I10").Value = 3.05
$ws.Range("K10").Value = 3.95
$ws.Range("L10").Value = 1.33
$ws.Range("P10").Value = 2.22
$ws.Range("Q10").Value = 1.7
$ws.Range("T10").Value = 1.6
$ws.Range("U10").Value = 2.38
$ws.Range("V10").Value = 1.48
$ws.Range("W10").Value = 1.62
$ws.Range("Z10").Value = 26
$ws.Range("AA10").Value = 55
$ws.Range("AD10").Value = 14
$ws.Range("AE10").Value = 32
$ws.Range("AF10").Value = 22
$ws.Range("AH10").Value = 16.5
$ws.Range("AI10").Value = 44
$ws.Range("AK10").Value = 26
$ws.Range("AM10").Value = 70
$ws.Range("AN10").Value = 17
$ws.Range("AO10").Value = 22
$ws.Range("H11").Value = 1.78
$ws.Range("N11").Value = 6.2
$ws.Range("S11").Value = 2.2
$ws.Range("U11").Value = 2.58
$ws.Range("AC11").Value = 13.5
$ws.Range("AD11").Value = 13.5
$ws.Range("F12").Value = 2.72
$ws.Range("G12").Value = 2.92
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 3.65
$ws.Range("L12").Value = 1.41
$ws.Range("O12").Value = 1.33
$ws.Range("Q12").Value = 1.99
$ws.Range("S12").Value = 3.55
$ws.Range("U12").Value = 1.93
$ws.Range("AH12").Value = 19.5
$ws.Range("F13").Value = 1.79
$ws.Range("G13").Value = 1.97
$ws.Range("H13").Value = 3.95
$ws.Range("T13").Value = 1.56
$ws.Range("V13").Value = 1.26
$ws.Range("AD13").Value = 22
$ws.Range("AH13").Value = 20
$ws.Range("AJ13").Value = 25
$ws.Range("F14").Value = 1.9
$ws.Range("G14").Value = 1.97
$ws.Range("I14").Value = 4.7
$ws.Range("K14").Value = 3.95
$ws.Range("T14").Value = 1.92
$ws.Range("U14").Value = 1.87
$ws.Range("V14").Value = 1.27
$ws.Range("X14").Value = 13.5
$ws.Range("Z14").Value = 40
$ws.Range("AC14").Value = 9.6
$ws.Range("AD14").Value = 22
$ws.Range("AE14").Value = 75
$ws.Range("AH14").Value = 26
$ws.Range("AN14").Value = 20
